# Auto-generated edit script: updates cryptos list Price (D) and Volume(1h) (E) columns
# to match scraped values from the commit "Updated cryptos list on Thu Mar  9 06:47:20 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): values that look like plain decimals would be auto-parsed as
# numbers by Excel, losing formatting (e.g. trailing zeros, "1.060" -> 1.06). Those cells
# are forced to Text format first so the literal string is preserved, matching the source
# data which stores these as inline text. Values already containing two "." separators
# (thousands-grouped, e.g. "21.740.42") have no valid numeric parse and stay text naturally,
# so we leave their cell format untouched.
$ws.Range("D2").Value = "21.740.42"
$ws.Range("D3").Value = "1.539.33"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.93"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3904"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3188"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.97"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07206"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.060"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.648"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.65"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.618"
$ws.Range("D16").Value = "1.541.19"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001110"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06588"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.26"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.157"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.41"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.87"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.404"
$ws.Range("D25").Value = "21.748.86"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.377"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.56"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.42"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.838"
$ws.Range("D30").Value = "1.715.93"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.70"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9769"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.925"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08200"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.833"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06094"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.142"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.482"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02203"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2040"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.188"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.70"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5755"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.743"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5514"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.164"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.29"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.873"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06717"

# --- Volume(1h) column (E): percentage strings with surrounding padding spaces; these
# already fail Excel's numeric/percentage auto-parse, so no format forcing is required.
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("E7").Value = "  +3.97%  "
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("E9").Value = "  +3.90%  "
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("E11").Value = "  -5.92%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("E14").Value = "  -4.47%  "
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("E21").Value = "  -4.45%  "
$ws.Range("E22").Value = "  -3.44%  "
$ws.Range("E23").Value = "  -5.90%  "
$ws.Range("E24").Value = "  +7.42%  "
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("E26").Value = "  -5.82%  "
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("E32").Value = "  -12.53%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("E35").Value = "  -4.71%  "
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("E38").Value = "  -7.94%  "
$ws.Range("E39").Value = "  -3.47%  "
$ws.Range("E40").Value = "  -3.21%  "
$ws.Range("E41").Value = "  -2.29%  "
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("E44").Value = "  -3.00%  "
$ws.Range("E45").Value = "  -3.73%  "
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("E47").Value = "  -3.88%  "
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("E49").Value = "  -2.00%  "
$ws.Range("E50").Value = "  -3.39%  "
$ws.Range("E51").Value = "  -3.02%  "
